$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows to append after the existing data (rows 2-4), reproducing
# rows 3, 4, then 2 again (as scraped a second time into the sheet).
$rows = @(
    @(" Dubai (DSC)", " October 04 2020", "Super Kings won by 10 wickets (with 14 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Sarfaraz Khan ", "14", "9", "2", "0", "155.55"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", "Sarfaraz Khan ", "7", "8", "1", "0", "87.50"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", "Sarfaraz Khan ", "12", "12", "2", "0", "100.00")
)

# Columns G:K hold numeric-looking values (runs, balls, 4s, 6s, strike rate)
# that must stay stored as text, same as the rest of the sheet.
$textCols = 7, 8, 9, 10, 11

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}

Write-Output "rows appended"
